# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from NMDP Rh Code Sys")

# 1. Rename the second worksheet.
$ws2.Name = "Include #0"

# 2. Bump the Version value (row 3, column B).
$ws1.Range("B3").Value = "0.1.1"

# 3. Update the Date value (row 8, column B).
$ws1.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# 4. Insert a new "Jurisdiction" property row right after "Contact" (row 10),
#    pushing Description / Purpose / Copyright / Immutable down by one row.
$ws1.Rows.Item(11).Insert()

# Copy the formatting of the row below (now the old "Description" row) onto
# the newly inserted blank row so the styling matches the rest of the table.
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

Write-Host "Applied NMDP Rh status codes metadata update"
